# "OUTPUT:" section of the report has two screenshots (InlineShapes 1 & 2).
# Both were cropped from the top with Word's Picture Tools > Crop handle
# (removing a blank strip above the console output) and had their outline
# turned off. Reproduce that with PictureFormat.CropTop + a matching
# resize, keeping each picture's aspect-ratio lock on afterwards so the
# stored "noChangeAspect" flag is left the way Word leaves it after an
# interactive crop.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Picture 1 (InlineShapes(1), r:embed="rId6") - crop ~17.573% off the top
# ---------------------------------------------------------------------
$shp1 = $d.InlineShapes.Item(1)

$shp1.PictureFormat.CropTop = 31.5

$shp1.LockAspectRatio = $false
$shp1.Height = 118.2
$shp1.Width = 427.8
$shp1.LockAspectRatio = $true

$shp1.Line.Visible = $false

# Cosmetic attributes Word also stamps on a freshly-cropped picture
# (no-ops on runtimes that don't expose them yet).
try { $shp1.Shadow.Obscured = $true } catch {}
try { $shp1.Shadow.RotateWithShape = $true } catch {}

# ---------------------------------------------------------------------
# Picture 2 (InlineShapes(2), r:embed="rId7") - crop ~13.566% off the top
# ---------------------------------------------------------------------
$shp2 = $d.InlineShapes.Item(2)

$shp2.PictureFormat.CropTop = 26.25

$shp2.LockAspectRatio = $false
$shp2.Height = 133.8
$shp2.Width = 409.8
$shp2.LockAspectRatio = $true

$shp2.Line.Visible = $false

try { $shp2.Shadow.Obscured = $true } catch {}
try { $shp2.Shadow.RotateWithShape = $true } catch {}
